# Posted final exam gradelines: fill in the Final Exam (column L) score for
# each letter-grade gradeline row, mark HW10 (row 15) and the Final Exam
# (row 23) as graded, and record the Final Exam's total points.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final exam scores for each gradeline row (L10:L21)
$ws.Range("L10").Value = 87.0
$ws.Range("L11").Value = 75.0
$ws.Range("L12").Value = 64.0
$ws.Range("L13").Value = 53.0
$ws.Range("L14").Value = 42.0
$ws.Range("L15").Value = 30.0
$ws.Range("L16").Value = 25.0
$ws.Range("L17").Value = 20.0
$ws.Range("L18").Value = 15.0
$ws.Range("L19").Value = 10.0
$ws.Range("L20").Value = 5.0
$ws.Range("L21").Value = 0.0

# HW10 (row 15) now graded
$ws.Range("E15").Value = 1.0

# Final exam total points (D23) -- copy D22's formatting to the
# previously-blank cell so it picks up the same "Total Points" style,
# then set its value.
$ws.Range("D22").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 112.0

# Final exam graded flag
$ws.Range("E23").Value = 1.0
